# Applies the "Update menu and exports" edit to the DaBeast Services Menu slide.
# Splits several plain runs into bold/italic-mixed runs (to emphasize certain
# words), fixes a couple of typos/price numbers, and merges the
# "Gatorade, Redbull, ..." runs back into a single run.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# ---------------------------------------------------------------------------
# Shape "Google Shape;94;p1" (Breakfast column)
# ---------------------------------------------------------------------------
$sh94 = $s.Shapes.Item(6)
$tr94 = $sh94.TextFrame.TextRange

# Paragraph 5: "Made from the freshest plantains"
#   -> "Made from the freshest " (unchanged) + "plantains" (bold)
$para = $tr94.Paragraphs(5)
$bold_part = $tr94.Characters($para.Start + 23, 9)   # "plantains"
$bold_part.Font.Bold = $true

# Paragraph 10: "\t    Fried cheese w/ guava marmalade"
#   -> "\t    Fried " (unchanged) + "cheese" (bold) + " w/ " (unchanged)
#      + "guava marmalade" (bold)
$para = $tr94.Paragraphs(10)
$base = $para.Start + 5                               # skip leading "\t    "
$tr94.Characters($base + 6, 6).Font.Bold = $true       # "cheese"
$tr94.Characters($base + 16, 15).Font.Bold = $true     # "guava marmalade"

# ---------------------------------------------------------------------------
# Shape "Google Shape;95;p1" (Lunch / Dinner column)
# ---------------------------------------------------------------------------
$sh95 = $s.Shapes.Item(7)
$tr95 = $sh95.TextFrame.TextRange

# Paragraph 3: "\t   Rice, beans w/ smoked chicken"
#   -> "Rice" (bold) + ", " + "beans" (bold) + " w/ " + "smoked chicken" (bold)
$para = $tr95.Paragraphs(3)
$base = $para.Start + 4                               # skip leading "\t   "
$tr95.Characters($base, 4).Font.Bold = $true           # "Rice"
$tr95.Characters($base + 6, 5).Font.Bold = $true       # "beans"
$tr95.Characters($base + 15, 14).Font.Bold = $true     # "smoked chicken"

# Paragraph 6: "\t   Smoked chicken skewers  "
#   -> whole dish name + trailing spaces become bold
$para = $tr95.Paragraphs(6)
$base = $para.Start + 4                               # skip leading "\t   "
$tr95.Characters($base, 22).Font.Bold = $true          # "Smoked chicken skewers"
$tr95.Characters($base + 22, 2).Font.Bold = $true      # trailing "  "

# Paragraph 10: "\t          Sweet plaintains" (typo) -> "... Sweet plantains" (bold)
$para = $tr95.Paragraphs(10)
$base = $para.Start + 11                              # skip leading "\t          "
$word = $tr95.Characters($base + 6, 10)               # "plaintains"
$word.Text = "plantains"
$word = $tr95.Characters($base + 6, 9)                # re-grab after shrink ("plantains")
$word.Font.Bold = $true

# Paragraph 12: "        Crisped plantains"
#   -> "        Crisped " (unchanged) + "plantains" (bold)
$para = $tr95.Paragraphs(12)
$tr95.Characters($para.Start + 16, 9).Font.Bold = $true  # "plantains"

# ---------------------------------------------------------------------------
# Shape "Google Shape;96;p1" (Desserts column)
# ---------------------------------------------------------------------------
$sh96 = $s.Shapes.Item(8)
$tr96 = $sh96.TextFrame.TextRange

# Paragraph 4: "Fried Cheesecake Bites …. 6" -> price raised to 7
$para = $tr96.Paragraphs(4)
$priceRun = $tr96.Characters($para.Start + 26, 1)     # "6"
$priceRun.Text = "7"

# ---------------------------------------------------------------------------
# Shape "Google Shape;99;p1" (Beverages column)
# ---------------------------------------------------------------------------
$sh99 = $s.Shapes.Item(11)
$tr99 = $sh99.TextFrame.TextRange

# Paragraph 6: "Gatorade, Redbull, Coke, Sprite, Fanta, Iced Tea, etc."
#   -> fix "Redbull" to "Red Bull" (merges back into a single plain run)
$para = $tr99.Paragraphs(6)
$word = $tr99.Characters($para.Start + 10, 7)          # "Redbull"
$word.Text = "Red Bull"
